$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows for the "Holden" and "Rizzie Spiral" measurement methods after row 3
$ws.Rows("4:5").Insert()

# Copy formatting (bold, border, centered) from A3 down to the new A4:A5 rank cells
$ws.Range("A3").Copy()
$ws.Range("A4:A5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New row 4: Holden
$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = "Holden"
$row4vals = @(1.030554051782185,0.9839906515144051,1.010632073398005,1.010632073398005,0.9558855347724781,1.005101250207399,0.9988463577875769,0.9288940756220727,1.010632073398005,0.9288940756220727,1.041987970143587,1.010632073398005,1.041987970143587,0.98544102288283,1.012989310828996,0.9938380397212215,0.9849575657600217,0.9938380397212215,0.9913761926695174,0.9952273688152149,0.9944864956534637)
for ($i = 0; $i -lt $row4vals.Length; $i++) {
    $ws.Cells.Item(4, 3 + $i).Value = $row4vals[$i]
}

# New row 5: Rizzie Spiral
$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(5, 2).Value = "Rizzie Spiral"
$row5vals = @(1.013615424979588,0.9787787558854281,1.185041148226768,1.185041148226768,0.9516066533321355,1.055722173225175,0.934497101450411,0.9565206442981622,1.185041148226768,0.9565206442981622,0.9348891019861739,1.185041148226768,0.9348891019861739,0.945704873142168,0.956833928935801,1.025483631503701,0.9567295007232547,1.025483631503701,1.013807412599133,1.04805415972466,1.00133387542298)
for ($i = 0; $i -lt $row5vals.Length; $i++) {
    $ws.Cells.Item(5, 3 + $i).Value = $row5vals[$i]
}

# Rename "Thomas Hex" -> "Matthies Hex" (now located at row 11 after the insert)
$ws.Range("B11").Value = "Matthies Hex"
